$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@("Shai Gilgeous-Alexander", "PG", "Oklahoma City Thunder")
    ,@("Keyonte George", "PG,SG", "Utah Jazz")
    ,@("Kyrie Irving", "PG,SG", "Dallas Mavericks")
    ,@("Dennis Schröder", "PG", "Brooklyn Nets")
    ,@("Lauri Markkanen", "SF,PF", "Utah Jazz")
    ,@("Tobias Harris", "SF,PF", "Detroit Pistons")
    ,@("Zach LaVine", "SG,SF", "Chicago Bulls")
    ,@("Christian Braun", "SG,SF", "Denver Nuggets")
    ,@("John Collins", "PF,C", "Utah Jazz")
    ,@("Dorian Finney-Smith", "SF,PF,C", "Brooklyn Nets")
    ,@("Jalen Williams", "SG,SF,PF", "Oklahoma City Thunder")
    ,@("Jordan Poole", "PG,SG", "Washington Wizards")
    ,@("Joel Embiid", "C", "Philadelphia 76ers")
    ,@("RJ Barrett", "SF,PF", "Toronto Raptors")
    ,@("Brandon Boston Jr.", "SG,SF", "New Orleans Pelicans")
    ,@("CJ McCollum", "PG,SG", "New Orleans Pelicans")
    ,@("Jimmy Butler", "SF,PF", "Miami Heat")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $data[$i][0]
    $ws.Cells.Item($r, 2).Value = $data[$i][1]
    $ws.Cells.Item($r, 3).Value = $data[$i][2]
}

# Row 19 no longer exists in the updated table; delete it entirely
$ws.Range("A19:C19").Delete()

